# Review Log Sheet update: "Add Tags of Test Cases to RTM"
#
# This adds a new reviewed entry (Rev_12_03) into the RTM review row that
# used to be the last (blank) row of the log (row 29), and pushes the
# previously-blank trailing row down to row 30, extending the sheet by one
# row overall (dimension C1:U37 -> C1:U38).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Row 30 keeps the same merged-cell layout every review row uses.
#        (Merging first, before the formatting is copied in below, keeps
#        Excel from "smart-splitting" the medium border across the newly
#        merged pair of cells.)
$ws.Range("C30:D30").Merge()
$ws.Range("E30:F30").Merge()
$ws.Range("G30:I30").Merge()
$ws.Range("J30:L30").Merge()
$ws.Range("M30:O30").Merge()
$ws.Range("P30:Q30").Merge()
$ws.Range("R30:S30").Merge()

# Duplicate the formatting of the (currently blank) row 29 onto the new
# blank row 30, so row 30 keeps looking like the old "spare" row.
$ws.Range("C29:U29").Copy()
$ws.Range("C30:U30").PasteSpecial(-4122)   # xlPasteFormats
# E30/F30 (Periority) lose the highlighted fill the active rows use, back to
# plain centered style -- same as C30/D30.
$ws.Range("C30").Copy()
$ws.Range("E30:F30").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(30).RowHeight = 15

# --- 2. Fill in the new review entry on row 29.
#        The Periority cell picks up the same highlighted fill used by the
#        other filled-in rows (copy its format from row 27's Periority cell).
$ws.Range("E27:F27").Copy()
$ws.Range("E29:F29").PasteSpecial(-4122)   # xlPasteFormats

#        (Comment cell is set first so the new shared strings land in the
#        same order as the source workbook: "Add Tags of Test Cases" before
#        "Rev_12_03".)
$ws.Range("J29").Value = "Add Tags of Test Cases"   # Comment
$ws.Range("C29").Value = "Rev_12_03"                 # Review No
$ws.Range("E29").Value = "Meduim"                    # Periority
$ws.Range("G29").Value = "RTM"                       # Document Name
$ws.Range("M29").Value = "Hossam"                    # Assigned Person
$ws.Range("P29").Value = (Get-Date -Year 2016 -Month 8 -Day 4 -Hour 0 -Minute 0 -Second 0)   # Start time
$ws.Range("R29").Value = (Get-Date -Year 2016 -Month 8 -Day 4 -Hour 0 -Minute 0 -Second 0)   # End time
$ws.Range("T29").Value = "DONE"                      # Status
$ws.Range("U29").Value = "Hossam"                    # Reviewer

# Row 29 grows slightly taller to match the other filled-in review rows.
$ws.Rows.Item(29).RowHeight = 15.6

# --- 3. A different review (row 26) shifted by a day: 8/4/2016 -> 8/5/2016.
$ws.Range("P26").Value = (Get-Date -Year 2016 -Month 8 -Day 5 -Hour 0 -Minute 0 -Second 0)

# --- 4. The sheet grew by one trailing blank row (row 38), matching the
#        look of the other blank rows (31-37) below the table.
$ws.Range("C37:U37").Copy()
$ws.Range("C38:U38").PasteSpecial(-4122)   # xlPasteFormats

# --- 5. Selection / scroll position, as left by the editing session.
$ws.Range("B15").Select()
$ws.Range("C29:D29").Select()
